# Auto-generated Excel COM-interop script applying scheduled-runner updates
# to the Halicarnassus_Profits leve-crafting-profit tables (per-sheet tables:
# ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Refreshed currentAveragePrice /
# LevePrice / LeveProfit columns (H, I, J, K, L, M, N) from the latest market
# data pull for the affected leve rows.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 13754.75
$ws.Range("I21").Value = 13754.75
$ws.Range("K21").Value = 13754.75
$ws.Range("M21").Value = -13286.75
$ws.Range("H23").Value = 13754.75
$ws.Range("I23").Value = 13754.75
$ws.Range("K23").Value = 13754.75
$ws.Range("M23").Value = -13520.75
$ws.Range("H33").Value = 107.36364
$ws.Range("I33").Value = 117.888885
$ws.Range("K33").Value = 117.888885
$ws.Range("M33").Value = 111.111115
$ws.Range("H51").Value = 1999.5
$ws.Range("I51").Value = 1400
$ws.Range("J51").Value = 2199.3333
$ws.Range("K51").Value = 1400
$ws.Range("L51").Value = 2199.3333
$ws.Range("M51").Value = -916
$ws.Range("N51").Value = -3167.3333
$ws.Range("H69").Value = 6999
$ws.Range("J69").Value = 6999
$ws.Range("L69").Value = 20997
$ws.Range("N69").Value = -22745
$ws.Range("H72").Value = 6999
$ws.Range("J72").Value = 6999
$ws.Range("L72").Value = 62991
$ws.Range("N72").Value = -71727
$ws.Range("H74").Value = 11391.667
$ws.Range("I74").Value = 4587.5
$ws.Range("K74").Value = 4587.5
$ws.Range("M74").Value = -3651.5
$ws.Range("H77").Value = 11391.667
$ws.Range("I77").Value = 4587.5
$ws.Range("K77").Value = 22937.5
$ws.Range("M77").Value = -18257.5
$ws.Range("H104").Value = 1200
$ws.Range("I104").Value = 1200
$ws.Range("K104").Value = 3600
$ws.Range("H121").Value = 877.5
$ws.Range("J121").Value = 877.5
$ws.Range("L121").Value = 2632.5
$ws.Range("N121").Value = -6126.5
$ws.Range("H127").Value = 2603
$ws.Range("I127").Value = 3200.8
$ws.Range("K127").Value = 9602.400000000001
$ws.Range("M127").Value = -4642.400000000001
$ws.Range("H137").Value = 2375.1924
$ws.Range("I137").Value = 1868.0667
$ws.Range("J137").Value = 3066.7273
$ws.Range("K137").Value = 5604.2001
$ws.Range("L137").Value = 9200.1819
$ws.Range("M137").Value = -3054.2001
$ws.Range("N137").Value = -14300.1819
$ws.Range("H141").Value = 1892
$ws.Range("I141").Value = 1491
$ws.Range("K141").Value = 4473
$ws.Range("M141").Value = 707
$ws.Range("M104").Value = -1853

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4425.3
$ws.Range("I32").Value = 4425.3
$ws.Range("K32").Value = 4425.3
$ws.Range("M32").Value = -4138.3
$ws.Range("H88").Value = 2054
$ws.Range("I88").Value = 1250
$ws.Range("J88").Value = 2456
$ws.Range("K88").Value = 1250
$ws.Range("L88").Value = 2456
$ws.Range("M88").Value = -844
$ws.Range("N88").Value = -3268
$ws.Range("H91").Value = 2054
$ws.Range("I91").Value = 1250
$ws.Range("J91").Value = 2456
$ws.Range("K91").Value = 1250
$ws.Range("L91").Value = 2456
$ws.Range("M91").Value = 154
$ws.Range("N91").Value = -5264
$ws.Range("H102").Value = 6491.154
$ws.Range("I102").Value = 5042.778
$ws.Range("K102").Value = 5042.778
$ws.Range("M102").Value = -3420.778

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 54750
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("H78").Value = 54750
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("H86").Value = 4911.5
$ws.Range("I86").Value = 455.85715
$ws.Range("J86").Value = 8377
$ws.Range("K86").Value = 455.85715
$ws.Range("L86").Value = 8377
$ws.Range("M86").Value = 667.14285
$ws.Range("N86").Value = -10623
$ws.Range("H89").Value = 4911.5
$ws.Range("I89").Value = 455.85715
$ws.Range("J89").Value = 8377
$ws.Range("K89").Value = 2279.28575
$ws.Range("L89").Value = 41885
$ws.Range("M89").Value = 3336.71425
$ws.Range("N89").Value = -53117
$ws.Range("H105").Value = 2010
$ws.Range("I105").Value = 2010
$ws.Range("K105").Value = 2010
$ws.Range("M105").Value = -263

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1739.8334
$ws.Range("I16").Value = 1747.5
$ws.Range("K16").Value = 1747.5
$ws.Range("M16").Value = -1460.5
$ws.Range("H31").Value = 5843.857
$ws.Range("I31").Value = 1372
$ws.Range("K31").Value = 1372
$ws.Range("M31").Value = -1077
$ws.Range("H34").Value = 5843.857
$ws.Range("I34").Value = 1372
$ws.Range("K34").Value = 1372
$ws.Range("M34").Value = -1170
$ws.Range("H35").Value = 55712.777
$ws.Range("I35").Value = 62630.625
$ws.Range("K35").Value = 62630.625
$ws.Range("M35").Value = -62336.625
$ws.Range("H99").Value = 3038.5833
$ws.Range("I99").Value = 2991.35
$ws.Range("K99").Value = 2991.35
$ws.Range("M99").Value = -1493.35
$ws.Range("H113").Value = 1739.8334
$ws.Range("I113").Value = 1747.5
$ws.Range("K113").Value = 1747.5
$ws.Range("M113").Value = 422.5
$ws.Range("H126").Value = 3038.5833
$ws.Range("I126").Value = 2991.35
$ws.Range("K126").Value = 8974.049999999999
$ws.Range("M126").Value = -6504.049999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 14925
$ws.Range("J43").Value = 14925
$ws.Range("L43").Value = 14925
$ws.Range("N43").Value = -15227
$ws.Range("H46").Value = 11987
$ws.Range("H80").Value = 2500
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 2500
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -24984

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 9177.888999999999
$ws.Range("I68").Value = 8125
$ws.Range("J68").Value = 10020.2
$ws.Range("K68").Value = 8125
$ws.Range("L68").Value = 10020.2
$ws.Range("M68").Value = -7376
$ws.Range("N68").Value = -11518.2
$ws.Range("H71").Value = 9177.888999999999
$ws.Range("I71").Value = 8125
$ws.Range("J71").Value = 10020.2
$ws.Range("K71").Value = 40625
$ws.Range("L71").Value = 50101
$ws.Range("M71").Value = -36881
$ws.Range("N71").Value = -57589
$ws.Range("H82").Value = 2545.2222
$ws.Range("I82").Value = 1174.25
$ws.Range("J82").Value = 3642
$ws.Range("K82").Value = 1174.25
$ws.Range("L82").Value = 3642
$ws.Range("M82").Value = -813.25
$ws.Range("N82").Value = -4364
$ws.Range("H85").Value = 2545.2222
$ws.Range("I85").Value = 1174.25
$ws.Range("J85").Value = 3642
$ws.Range("K85").Value = 1174.25
$ws.Range("L85").Value = 3642
$ws.Range("M85").Value = 73.75
$ws.Range("N85").Value = -6138
$ws.Range("H122").Value = 3997
$ws.Range("I122").Value = 3829.8333
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 11489.4999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -9039.499899999999
$ws.Range("N122").Value = -19900

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8055.4443
$ws.Range("I62").Value = 5875
$ws.Range("J62").Value = 9799.799999999999
$ws.Range("K62").Value = 5875
$ws.Range("L62").Value = 9799.799999999999
$ws.Range("M62").Value = -5251
$ws.Range("N62").Value = -11047.8
$ws.Range("H65").Value = 8055.4443
$ws.Range("I65").Value = 5875
$ws.Range("J65").Value = 9799.799999999999
$ws.Range("K65").Value = 29375
$ws.Range("L65").Value = 48999
$ws.Range("M65").Value = -26255
$ws.Range("N65").Value = -55239
$ws.Range("H132").Value = 2847.7273
$ws.Range("I132").Value = 2847.7273
$ws.Range("K132").Value = 8543.1819
$ws.Range("M132").Value = -6013.1819

# ---- Sheet: BSM (cell removals) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M75").ClearContents()
$ws.Range("M78").ClearContents()

Write-Host "Applied scheduled market-data refresh to 8 leve-profit sheets."